# Applies the "corrected reading of make financial results" update:
#  - shifts the simulation StartTime/StopTime window forward by 5 years
#  - refreshes scenario_data_emlab fuel/CO2 price assumptions for the new year
#  - replaces the conventionals / renewables / biogas unit tables with the
#    corrected (deduplicated) data read from the traderes DB, and adds the
#    missing storage unit row

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# times: StartTime / StopTime
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("times")
$ws.Range("B2").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B3").NumberFormat = "yyyy-mm-dd h:mm:ss"
$ws.Range("B2").Value = 45657.99861111111
$ws.Range("B3").Value = 46021.99861111111

# ---------------------------------------------------------------------
# scenario_data_emlab: base year + fuel / CO2 prices
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("scenario_data_emlab")
$ws.Range("B1").Value = 2025
$ws.Range("B2").Value = 37
$ws.Range("B5").Value = 10.895
$ws.Range("B6").Value = 23.92333333333333
$ws.Range("B7").Value = 56.19333333333333

# ---------------------------------------------------------------------
# conventionals: drop the stale duplicate rows, keep two corrected ones
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("conventionals")
$ws.Range("A4:A7").EntireRow.Delete()

$ws.Range("B2").Value = 99991700006
$ws.Range("C2").Value = "NATURAL_GAS"
$ws.Range("D2").Value = 4.5
$ws.Range("E2").Value = 0.43
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1

$ws.Range("B3").Value = 99990300008
$ws.Range("C3").Value = "NATURAL_GAS"
$ws.Range("D3").Value = 4.2
$ws.Range("E3").Value = 0.61
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1

# ---------------------------------------------------------------------
# renewables: drop the stale last row, renumber / correct the remaining 3
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("renewables")
$ws.Range("A5").EntireRow.Delete()

$ws.Range("B2").Value = 99992100002
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = "OtherPV"

$ws.Range("B3").Value = 99992400003
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 1.35
$ws.Range("E3").Value = "WindOn"

$ws.Range("B4").Value = 99992300007
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 2.7
$ws.Range("E4").Value = "WindOff"

# ---------------------------------------------------------------------
# storages: add the missing STORAGE unit row
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("storages")
$ws.Range("A2").Value = 0
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160
$ws.Range("A2").Borders.LineStyle = 1
$ws.Range("B2").Value = 99992600009
$ws.Range("C2").Value = "STORAGE"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 1

# ---------------------------------------------------------------------
# biogas: corrected identifier / installed power
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("biogas")
$ws.Range("B2").Value = 99990100004
$ws.Range("C2").Value = 1
